$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 'Digit_before_after'
$ws.Range("H3").Value = 'DIGIT_POST'
$ws.Range("K3").Value = 8
$ws.Range("H4").Value = 'Digit_before_after'
$ws.Range("K4").Value = 2
$ws.Range("H5").Value = 'DIGIT_POST'
$ws.Range("K5").Value = 7
$ws.Range("H6").Value = 'DICHOTIC_PRE'
$ws.Range("K6").Value = 8
$ws.Range("H7").Value = 'DICHOTIC_POST'
$ws.Range("K8").Value = 3
$ws.Range("H9").Value = 'DIGIT_PRE'
$ws.Range("H10").Value = 'DICHOTIC_PRE'
$ws.Range("K10").Value = 6
$ws.Range("H11").Value = 'DICHOTIC_POST'
$ws.Range("K11").Value = 8
$ws.Range("H12").Value = 'DICHOTIC_PRE'
$ws.Range("K12").Value = 5
$ws.Range("H14").Value = 'DIGIT_PRE'
$ws.Range("K14").Value = 3
$ws.Range("H15").Value = 'DIGIT_PRE'
$ws.Range("K16").Value = 4
$ws.Range("H17").Value = 'Dichotic_before_after'
$ws.Range("K17").Value = 8
$ws.Range("K18").Value = 4
$ws.Range("H19").Value = 'DICHOTIC_POST'
$ws.Range("K19").Value = 1
$ws.Range("K20").Value = 6
$ws.Range("H21").Value = 'Digit_before_after'
$ws.Range("K21").Value = 3
$ws.Range("H22").Value = 'DIGIT_PRE'
$ws.Range("H23").Value = 'DICHOTIC_PRE'
$ws.Range("K23").Value = 4
$ws.Range("K24").Value = 7
$ws.Range("H25").Value = 'DIGIT_POST'
$ws.Range("K25").Value = 8
$ws.Range("H26").Value = 'DICHOTIC_POST'
$ws.Range("K26").Value = 1
$ws.Range("H27").Value = 'Dichotic_before_after'
$ws.Range("K27").Value = 7
$ws.Range("H28").Value = 'DICHOTIC_POST'
$ws.Range("K28").Value = 3
$ws.Range("H29").Value = 'DICHOTIC_PRE'
$ws.Range("K29").Value = 7
$ws.Range("H30").Value = 'DIGIT_POST'
$ws.Range("K30").Value = 8
$ws.Range("H31").Value = 'Dichotic_before_after'
$ws.Range("K31").Value = 5
$ws.Range("K32").Value = 2
$ws.Range("H33").Value = 'Digit_before_after'
$ws.Range("K33").Value = 6
$ws.Range("H34").Value = 'DIGIT_PRE'
$ws.Range("K34").Value = 5
$ws.Range("H35").Value = 'DIGIT_POST'
$ws.Range("K35").Value = 6
$ws.Range("K36").Value = 6
$ws.Range("H37").Value = 'Digit_before_after'
$ws.Range("K37").Value = 5
$ws.Range("H38").Value = 'DIGIT_PRE'
$ws.Range("K38").Value = 5
$ws.Range("H39").Value = 'Dichotic_before_after'
$ws.Range("K39").Value = 7
$ws.Range("K40").Value = 7
$ws.Range("H41").Value = 'DIGIT_POST'
$ws.Range("K41").Value = 3
$ws.Range("H42").Value = 'Digit_before_after'
$ws.Range("H43").Value = 'DIGIT_POST'
$ws.Range("K43").Value = 7
$ws.Range("H44").Value = 'DICHOTIC_POST'
$ws.Range("K44").Value = 8
$ws.Range("H45").Value = 'DICHOTIC_PRE'
$ws.Range("K45").Value = 4
$ws.Range("H46").Value = 'DIGIT_POST'
$ws.Range("K46").Value = 2
$ws.Range("H47").Value = 'DIGIT_PRE'
$ws.Range("K47").Value = 4
$ws.Range("H48").Value = 'DIGIT_PRE'
$ws.Range("K48").Value = 7
$ws.Range("H49").Value = 'DIGIT_PRE'
$ws.Range("K49").Value = 8
$ws.Range("H50").Value = 'DIGIT_PRE'
$ws.Range("K50").Value = 6
$ws.Range("H51").Value = 'DICHOTIC_PRE'
$ws.Range("K51").Value = 3
$ws.Range("H52").Value = 'Digit_before_after'
$ws.Range("K52").Value = 8
$ws.Range("H53").Value = 'DIGIT_POST'
$ws.Range("K53").Value = 7
$ws.Range("K54").Value = 5
$ws.Range("K55").Value = 2
$ws.Range("K56").Value = 1
$ws.Range("H57").Value = 'DICHOTIC_POST'
$ws.Range("K57").Value = 8
$ws.Range("H58").Value = 'DIGIT_POST'
$ws.Range("K58").Value = 5
$ws.Range("H59").Value = 'DIGIT_PRE'
$ws.Range("K59").Value = 3
$ws.Range("H60").Value = 'Digit_before_after'
$ws.Range("K60").Value = 7
$ws.Range("H61").Value = 'Dichotic_before_after'
$ws.Range("K61").Value = 2
$ws.Range("H62").Value = 'DICHOTIC_POST'
$ws.Range("K62").Value = 2
$ws.Range("H63").Value = 'DICHOTIC_POST'
$ws.Range("K63").Value = 5
$ws.Range("H64").Value = 'Dichotic_before_after'
$ws.Range("K64").Value = 4
$ws.Range("H65").Value = 'DICHOTIC_POST'
$ws.Range("K65").Value = 6
$ws.Range("H66").Value = 'DICHOTIC_PRE'
$ws.Range("K66").Value = 5
$ws.Range("K68").Value = 8
$ws.Range("H69").Value = 'Digit_before_after'
$ws.Range("K69").Value = 2
$ws.Range("H70").Value = 'Dichotic_before_after'
$ws.Range("H71").Value = 'DICHOTIC_PRE'
$ws.Range("K71").Value = 1
$ws.Range("H72").Value = 'DIGIT_POST'
$ws.Range("K72").Value = 1
$ws.Range("H73").Value = 'Dichotic_before_after'
$ws.Range("K73").Value = 2
$ws.Range("K74").Value = 8
$ws.Range("H75").Value = 'DIGIT_PRE'
$ws.Range("K75").Value = 7
$ws.Range("H76").Value = 'Dichotic_before_after'
$ws.Range("K76").Value = 2
$ws.Range("H77").Value = 'DICHOTIC_PRE'
$ws.Range("K77").Value = 4
$ws.Range("H78").Value = 'Digit_before_after'
$ws.Range("K78").Value = 4
$ws.Range("H79").Value = 'Dichotic_before_after'
$ws.Range("K79").Value = 7
$ws.Range("K80").Value = 5
$ws.Range("H81").Value = 'DICHOTIC_POST'
